$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 600
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 600
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -826

$ws.Range("H17").Value = 2994.6667
$ws.Range("I17").Value = 800
$ws.Range("K17").Value = 2400
$ws.Range("M17").Value = -2232

$ws.Range("H28").Value = 103.75
$ws.Range("J28").Value = 72.5
$ws.Range("L28").Value = 72.5
$ws.Range("N28").Value = -1042.5

$ws.Range("H40").Value = 7766.6665
$ws.Range("I40").Value = 7766.6665
$ws.Range("K40").Value = 7766.6665
$ws.Range("M40").Value = -7591.6665

$ws.Range("H41").Value = 476.42856
$ws.Range("I41").Value = 261.66666
$ws.Range("K41").Value = 261.66666
$ws.Range("M41").Value = 178.33334

$ws.Range("H53").Value = 266.73334
$ws.Range("I53").Value = 243.625
$ws.Range("K53").Value = 243.625
$ws.Range("M53").Value = 393.375

$ws.Range("H58").Value = 1571.25
$ws.Range("I58").Value = 1571.25
$ws.Range("K58").Value = 4713.75
$ws.Range("M58").Value = -4563.75

$ws.Range("H69").Value = 260000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 260000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 780000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -781748

$ws.Range("H72").Value = 260000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 260000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 2340000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -2348736

$ws.Range("H82").Value = 4625
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 4625
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H92").Value = 830.44446
$ws.Range("I92").Value = 796.875
$ws.Range("K92").Value = 796.875
$ws.Range("M92").Value = 451.125

$ws.Range("H107").Value = 940.2143
$ws.Range("I107").Value = 1147.1818
$ws.Range("J107").Value = 181.33333
$ws.Range("K107").Value = 1147.1818
$ws.Range("L107").Value = 181.33333
$ws.Range("M107").Value = 772.8181999999999
$ws.Range("N107").Value = -4021.33333

$ws.Range("H137").Value = 2275.2424
$ws.Range("I137").Value = 2340.0417
$ws.Range("K137").Value = 7020.125100000001
$ws.Range("M137").Value = -4470.125100000001

$ws.Range("H138").Value = 4187.1
$ws.Range("I138").Value = 898.8570999999999
$ws.Range("J138").Value = 11859.667
$ws.Range("K138").Value = 2696.5713
$ws.Range("L138").Value = 35579.001
$ws.Range("M138").Value = 2443.4287
$ws.Range("N138").Value = -45859.001

$ws.Range("H141").Value = 850
$ws.Range("I141").Value = 850
$ws.Range("K141").Value = 2550
$ws.Range("M141").Value = 2630

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3851.5
$ws.Range("I2").Value = 5202.25
$ws.Range("J2").Value = 1150
$ws.Range("K2").Value = 5202.25
$ws.Range("L2").Value = 1150
$ws.Range("M2").Value = -5089.25
$ws.Range("N2").Value = -1376

$ws.Range("H61").Value = 1733.3334
$ws.Range("I61").Value = 1733.3334
$ws.Range("K61").Value = 1733.3334
$ws.Range("M61").Value = -1521.3334

$ws.Range("H74").Value = 11622.875
$ws.Range("I74").Value = 11885.05
$ws.Range("K74").Value = 11885.05
$ws.Range("M74").Value = -11011.05

$ws.Range("H77").Value = 11622.875
$ws.Range("I77").Value = 11885.05
$ws.Range("K77").Value = 59425.25
$ws.Range("M77").Value = -55057.25

$ws.Range("H116").Value = 3851.5
$ws.Range("I116").Value = 5202.25
$ws.Range("J116").Value = 1150
$ws.Range("K116").Value = 5202.25
$ws.Range("L116").Value = 1150
$ws.Range("M116").Value = -2908.25
$ws.Range("N116").Value = -5738

$ws.Range("H136").Value = 1733.3334
$ws.Range("I136").Value = 1733.3334
$ws.Range("K136").Value = 5200.0002
$ws.Range("M136").Value = -2650.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3851.5
$ws.Range("I3").Value = 5202.25
$ws.Range("J3").Value = 1150
$ws.Range("K3").Value = 5202.25
$ws.Range("L3").Value = 1150
$ws.Range("M3").Value = -5088.25
$ws.Range("N3").Value = -1378

$ws.Range("H134").Value = 2732.2632
$ws.Range("I134").Value = 2647.8235
$ws.Range("J134").Value = 3450
$ws.Range("K134").Value = 7943.470499999999
$ws.Range("L134").Value = 10350
$ws.Range("M134").Value = -5408.470499999999
$ws.Range("N134").Value = -15420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 26999.334
$ws.Range("I8").Value = 26998
$ws.Range("K8").Value = 26998
$ws.Range("M8").Value = -26858

$ws.Range("H22").Value = 308
$ws.Range("I22").Value = 309.6
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 309.6
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 40.39999999999998
$ws.Range("N22").Value = -1000

$ws.Range("H31").Value = 1729.7826
$ws.Range("J31").Value = 1882.1111
$ws.Range("L31").Value = 1882.1111
$ws.Range("N31").Value = -2472.1111

$ws.Range("H34").Value = 1729.7826
$ws.Range("J34").Value = 1882.1111
$ws.Range("L34").Value = 1882.1111
$ws.Range("N34").Value = -2286.1111

$ws.Range("H132").Value = 3358.353
$ws.Range("I132").Value = 3166.3333
$ws.Range("K132").Value = 9498.999899999999
$ws.Range("M132").Value = -6968.999899999999

$ws.Range("H134").Value = 4304.7144
$ws.Range("I134").Value = 4043.0833
$ws.Range("J134").Value = 5874.5
$ws.Range("K134").Value = 12129.2499
$ws.Range("L134").Value = 17623.5
$ws.Range("M134").Value = -9594.249899999999
$ws.Range("N134").Value = -22693.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3835.1428
$ws.Range("J75").Value = 2974.3333
$ws.Range("L75").Value = 8922.999899999999
$ws.Range("N75").Value = -10918.9999

$ws.Range("H78").Value = 3835.1428
$ws.Range("J78").Value = 2974.3333
$ws.Range("L78").Value = 26768.9997
$ws.Range("N78").Value = -36752.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5500
$ws.Range("I97").Value = 5500
$ws.Range("K97").Value = 5500
$ws.Range("M97").Value = -5004

$ws.Range("H107").Value = 881.1667
$ws.Range("J107").Value = 397.33334
$ws.Range("L107").Value = 397.33334
$ws.Range("N107").Value = -4237.33334

$ws.Range("H132").Value = 3158.1875
$ws.Range("J132").Value = 3999.5
$ws.Range("L132").Value = 11998.5
$ws.Range("N132").Value = -17058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13042.857
$ws.Range("I22").Value = 17900
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 17900
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -17605
$ws.Range("N22").Value = -1490

$ws.Range("H27").Value = 13042.857
$ws.Range("I27").Value = 17900
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 17900
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -17793
$ws.Range("N27").Value = -1114

$ws.Range("H30").Value = 838.6667
$ws.Range("I30").Value = 838.6667
$ws.Range("K30").Value = 838.6667
$ws.Range("M30").Value = -730.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 210323
$ws.Range("I2").Value = 210323
$ws.Range("K2").Value = 210323
$ws.Range("M2").Value = -210211

$ws.Range("H62").Value = 7692.857
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376

$ws.Range("H65").Value = 7692.857
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880

$ws.Range("H107").Value = 1328.4286
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 1259.8
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 3779.4
$ws.Range("M107").Value = -2580
$ws.Range("N107").Value = -7619.4

$ws.Range("H113").Value = 749.6
$ws.Range("I113").Value = 805.1111
$ws.Range("J113").Value = 250
$ws.Range("K113").Value = 2415.3333
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = -245.3332999999998
$ws.Range("N113").Value = -5090

$ws.Range("H132").Value = 722.5714
$ws.Range("I132").Value = 762.7692
$ws.Range("K132").Value = 2288.3076
$ws.Range("M132").Value = 241.6923999999999
